# Fix the placeholder "xxx" hotel-name values on Sheet1 with real hotel
# names, highlight the first one in bold, and set the sheet's print/page
# setup (A4 portrait).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2 ("xxx" -> "Grand Plaza Apartments") gets bolded.
$ws.Range("B2").Value = "Grand Plaza Apartments"
$ws.Range("B2").Font.Bold = $true

# B3 ("xxx" -> "Jumeirah Beach Hotel") stays regular weight.
$ws.Range("B3").Value = "Jumeirah Beach Hotel"

# Page setup: A4 paper, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
